$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 351
$hyperlinkCols = @(19, 20, 22, 23, 24, 25)  # S, T, V, W, X, Y

for ($r = 2; $r -le $lastRow; $r++) {
    # Update "Förändrad" (changed) date in column C from 45184 -> 45186
    $ws.Cells.Item($r, 3).Value = 45186

    $idCell = $ws.Cells.Item($r, 1)
    $idValue = $idCell.Value()

    if ($idValue -ne $null -and $idValue -ne "") {
        foreach ($col in $hyperlinkCols) {
            $cell = $ws.Cells.Item($r, $col)
            $f = $cell.Formula()
            if ($f -ne $null -and $f -ne "" -and $f.ToUpper().Contains("HYPERLINK") -and -not $f.Contains(",")) {
                $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $idValue + '")'
                $cell.Formula = $newFormula
            }
        }
    }
}
